$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "680÷9=75, 5" "121÷4=30, 1"
Replace-Text "233÷3=77, 2" "384÷5=76, 4"
Replace-Text "305÷4=76, 1" "765÷6=127, 3"
Replace-Text "867÷4=216, 3" "139÷7=19, 6"
Replace-Text "473÷7=67, 4" "809÷6=134, 5"
Replace-Text "610÷2=305, 0" "756÷8=94, 4"
Replace-Text "763÷2=381, 1" "607÷8=75, 7"
Replace-Text "402÷7=57, 3" "307÷2=153, 1"
Replace-Text "522÷7=74, 4" "990÷3=330, 0"
Replace-Text "978÷8=122, 2" "882÷9=98, 0"
Replace-Text "270÷6=45, 0" "446÷5=89, 1"
Replace-Text "173÷8=21, 5" "428÷2=214, 0"
Replace-Text "389÷6=64, 5" "928÷9=103, 1"
Replace-Text "311÷2=155, 1" "465÷3=155, 0"
Replace-Text "624÷7=89, 1" "455÷9=50, 5"
Replace-Text "623÷5=124, 3" "199÷2=99, 1"
Replace-Text "427÷7=61, 0" "485÷7=69, 2"
Replace-Text "513÷4=128, 1" "706÷6=117, 4"
Replace-Text "346÷7=49, 3" "521÷4=130, 1"
Replace-Text "577÷5=115, 2" "351÷5=70, 1"
Replace-Text "127÷2=63, 1" "926÷9=102, 8"
Replace-Text "169÷4=42, 1" "487÷3=162, 1"
Replace-Text "781÷7=111, 4" "264÷6=44, 0"
Replace-Text "762÷3=254, 0" "429÷4=107, 1"
Replace-Text "673÷6=112, 1" "193÷7=27, 4"

Write-Host "Done applying replacements"
